$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1888.826
$ws.Range("I40").Value = 1721
$ws.Range("J40").Value = 2314.8462
$ws.Range("K40").Value = 1721
$ws.Range("L40").Value = 2314.8462
$ws.Range("M40").Value = -1546
$ws.Range("N40").Value = -2664.8462
$ws.Range("H41").Value = 476.42856
$ws.Range("I41").Value = 437.8
$ws.Range("J41").Value = 573
$ws.Range("K41").Value = 437.8
$ws.Range("L41").Value = 573
$ws.Range("M41").Value = 2.199999999999989
$ws.Range("N41").Value = -1453
$ws.Range("H80").Value = 460.41666
$ws.Range("J80").Value = 435.77777
$ws.Range("L80").Value = 1307.33331
$ws.Range("N80").Value = -3303.33331
$ws.Range("H83").Value = 460.41666
$ws.Range("J83").Value = 435.77777
$ws.Range("L83").Value = 3921.99993
$ws.Range("N83").Value = -13905.99993
$ws.Range("H104").Value = 100
$ws.Range("I104").Value = 100
$ws.Range("K104").Value = 300
$ws.Range("M104").Value = 1447
$ws.Range("H111").Value = 44109.555
$ws.Range("I111").Value = 43833.168
$ws.Range("K111").Value = 131499.504
$ws.Range("M111").Value = -128432.504
$ws.Range("H116").Value = 11300.333
$ws.Range("I116").Value = 15001
$ws.Range("J116").Value = 7599.6665
$ws.Range("K116").Value = 15001
$ws.Range("L116").Value = 7599.6665
$ws.Range("M116").Value = -11559
$ws.Range("N116").Value = -14483.6665
$ws.Range("H132").Value = 2902.0908
$ws.Range("I132").Value = 2448.111
$ws.Range("J132").Value = 4945
$ws.Range("K132").Value = 7344.333
$ws.Range("L132").Value = 14835
$ws.Range("M132").Value = -4814.333
$ws.Range("N132").Value = -19895
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1750.5
$ws.Range("I2").Value = 1667.3334
$ws.Range("J2").Value = 2000
$ws.Range("K2").Value = 1667.3334
$ws.Range("L2").Value = 2000
$ws.Range("M2").Value = -1554.3334
$ws.Range("N2").Value = -2226
$ws.Range("H4").Value = 816.1667
$ws.Range("I4").Value = 474.5
$ws.Range("K4").Value = 474.5
$ws.Range("M4").Value = -358.5
$ws.Range("H116").Value = 1750.5
$ws.Range("I116").Value = 1667.3334
$ws.Range("J116").Value = 2000
$ws.Range("K116").Value = 1667.3334
$ws.Range("L116").Value = 2000
$ws.Range("M116").Value = 626.6666
$ws.Range("N116").Value = -6588
$ws.Range("H122").Value = 2996.6875
$ws.Range("I122").Value = 3124.6667
$ws.Range("J122").Value = 2919.9
$ws.Range("K122").Value = 9374.000100000001
$ws.Range("L122").Value = 8759.700000000001
$ws.Range("M122").Value = -6924.000100000001
$ws.Range("N122").Value = -13659.7
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1750.5
$ws.Range("I3").Value = 1667.3334
$ws.Range("J3").Value = 2000
$ws.Range("K3").Value = 1667.3334
$ws.Range("L3").Value = 2000
$ws.Range("M3").Value = -1553.3334
$ws.Range("N3").Value = -2228
$ws.Range("H88").Value = 21277.4
$ws.Range("I88").Value = 8497
$ws.Range("J88").Value = 24472.5
$ws.Range("K88").Value = 8497
$ws.Range("L88").Value = 24472.5
$ws.Range("M88").Value = -8091
$ws.Range("N88").Value = -25284.5
$ws.Range("H91").Value = 21277.4
$ws.Range("I91").Value = 8497
$ws.Range("J91").Value = 24472.5
$ws.Range("K91").Value = 8497
$ws.Range("L91").Value = 24472.5
$ws.Range("M91").Value = -7093
$ws.Range("N91").Value = -27280.5
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 810.9
$ws.Range("J5").Value = 1173
$ws.Range("L5").Value = 1173
$ws.Range("N5").Value = -1397
$ws.Range("H31").Value = 2738.0833
$ws.Range("I31").Value = 2647.5557
$ws.Range("J31").Value = 3009.6667
$ws.Range("K31").Value = 2647.5557
$ws.Range("L31").Value = 3009.6667
$ws.Range("M31").Value = -2352.5557
$ws.Range("N31").Value = -3599.6667
$ws.Range("H34").Value = 2738.0833
$ws.Range("I34").Value = 2647.5557
$ws.Range("J34").Value = 3009.6667
$ws.Range("K34").Value = 2647.5557
$ws.Range("L34").Value = 3009.6667
$ws.Range("M34").Value = -2445.5557
$ws.Range("N34").Value = -3413.6667
$ws.Range("H134").Value = 3037.818
$ws.Range("I134").Value = 2813.7856
$ws.Range("J134").Value = 3429.875
$ws.Range("K134").Value = 8441.356800000001
$ws.Range("L134").Value = 10289.625
$ws.Range("M134").Value = -5906.356800000001
$ws.Range("N134").Value = -15359.625
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 141.25
$ws.Range("I11").Value = 125
$ws.Range("J11").Value = 146.66667
$ws.Range("K11").Value = 375
$ws.Range("L11").Value = 440.00001
$ws.Range("M11").Value = -235
$ws.Range("N11").Value = -720.00001
$ws.Range("H37").Value = 69944
$ws.Range("J37").Value = 69944
$ws.Range("L37").Value = 209832
$ws.Range("N37").Value = -210056
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 55.666668
$ws.Range("I2").Value = 34
$ws.Range("J2").Value = 120.666664
$ws.Range("K2").Value = 34
$ws.Range("L2").Value = 120.666664
$ws.Range("M2").Value = 79
$ws.Range("N2").Value = -346.666664
$ws.Range("H15").Value = 54880
$ws.Range("J15").Value = 54880
$ws.Range("L15").Value = 54880
$ws.Range("N15").Value = -55456
$ws.Range("H28").Value = 10000
$ws.Range("J28").Value = 10000
$ws.Range("L28").Value = 10000
$ws.Range("N28").Value = -10384
$ws.Range("H69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("H72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("H81").Value = 54880
$ws.Range("J81").Value = 54880
$ws.Range("L81").Value = 54880
$ws.Range("N81").Value = -56876
$ws.Range("H84").Value = 54880
$ws.Range("J84").Value = 54880
$ws.Range("L84").Value = 164640
$ws.Range("N84").Value = -174624
$ws.Range("H94").Value = 26706.285
$ws.Range("I94").Value = 17740.334
$ws.Range("J94").Value = 29151.545
$ws.Range("K94").Value = 17740.334
$ws.Range("L94").Value = 29151.545
$ws.Range("M94").Value = -17064.334
$ws.Range("N94").Value = -30503.545
$ws.Range("H122").Value = 8308.875
$ws.Range("I122").Value = 7394.6
$ws.Range("K122").Value = 22183.8
$ws.Range("M122").Value = -19733.8
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H30").Value = 1014.5
$ws.Range("I30").Value = 1014.5
$ws.Range("K30").Value = 1014.5
$ws.Range("M30").Value = -906.5
$ws.Range("H46").Value = 1804.8889
$ws.Range("I46").Value = 1449.6666
$ws.Range("K46").Value = 1449.6666
$ws.Range("M46").Value = -1261.6666
$ws.Range("H55").Value = 1413.4286
$ws.Range("I55").Value = 599
$ws.Range("K55").Value = 599
$ws.Range("M55").Value = -426
$ws.Range("H93").Value = 1271.909
$ws.Range("I93").Value = 1070.8572
$ws.Range("J93").Value = 1623.75
$ws.Range("K93").Value = 1070.8572
$ws.Range("L93").Value = 1623.75
$ws.Range("M93").Value = 177.1428000000001
$ws.Range("N93").Value = -4119.75
$ws.Range("H122").Value = 5830.7144
$ws.Range("I122").Value = 4657.5
$ws.Range("K122").Value = 13972.5
$ws.Range("M122").Value = -11522.5
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 2705.35
$ws.Range("I107").Value = 2735.0588
$ws.Range("K107").Value = 8205.1764
$ws.Range("M107").Value = -6285.1764
$ws.Range("H122").Value = 1477.7778
$ws.Range("I122").Value = 1487.5
$ws.Range("K122").Value = 4462.5
$ws.Range("M122").Value = -2012.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("N69").ClearContents()
$ws.Range("N72").ClearContents()
